# Edit script: add "metadata" sheet and refresh "time_taken" / query timestamps
# on the "data" sheet, per the commit "Refined metadata to be additional tab".

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Update the per-row query timestamps in column F of the "data" sheet.
# ---------------------------------------------------------------------------
$timeTaken = @{
    2 = "2021-10-05 14:35:54.275030"
    3 = "2021-10-05 14:35:54.275037"
    4 = "2021-10-05 14:35:54.275040"
    5 = "2021-10-05 14:35:54.275042"
    6 = "2021-10-05 14:35:54.275045"
    7 = "2021-10-05 14:35:54.275047"
    8 = "2021-10-05 14:35:54.275049"
    9 = "2021-10-05 14:35:54.275052"
    10 = "2021-10-05 14:35:54.275054"
    11 = "2021-10-05 14:35:54.275057"
    12 = "2021-10-05 14:35:54.275059"
    13 = "2021-10-05 14:35:54.275061"
    14 = "2021-10-05 14:35:54.275064"
    15 = "2021-10-05 14:35:54.275066"
    16 = "2021-10-05 14:35:54.275068"
    17 = "2021-10-05 14:35:54.275071"
    18 = "2021-10-05 14:35:54.275074"
    19 = "2021-10-05 14:35:54.275076"
    20 = "2021-10-05 14:35:54.275079"
    21 = "2021-10-05 14:35:54.275081"
    22 = "2021-10-05 14:35:54.275083"
    23 = "2021-10-05 14:35:54.275086"
    24 = "2021-10-05 14:35:54.275088"
    25 = "2021-10-05 14:35:54.275090"
    26 = "2021-10-05 14:35:54.275093"
    27 = "2021-10-05 14:35:54.275095"
    28 = "2021-10-05 14:35:54.275097"
    29 = "2021-10-05 14:35:54.275100"
    30 = "2021-10-05 14:35:54.275102"
    31 = "2021-10-05 14:35:54.275104"
    32 = "2021-10-05 14:35:54.275107"
    33 = "2021-10-05 14:35:54.275109"
    34 = "2021-10-05 14:35:54.275111"
    35 = "2021-10-05 14:35:54.275114"
    36 = "2021-10-05 14:35:54.275116"
    37 = "2021-10-05 14:35:54.275118"
    38 = "2021-10-05 14:35:54.275121"
    39 = "2021-10-05 14:35:54.275123"
    40 = "2021-10-05 14:35:54.275125"
    41 = "2021-10-05 14:35:54.275127"
    42 = "2021-10-05 14:35:54.275130"
    43 = "2021-10-05 14:35:54.275133"
    44 = "2021-10-05 14:35:54.275135"
    45 = "2021-10-05 14:35:54.275137"
    46 = "2021-10-05 14:35:54.275139"
    47 = "2021-10-05 14:35:54.275142"
    48 = "2021-10-05 14:35:54.275144"
    49 = "2021-10-05 14:35:54.275146"
    50 = "2021-10-05 14:35:54.275149"
    51 = "2021-10-05 14:35:54.275151"
    52 = "2021-10-05 14:35:54.275153"
    53 = "2021-10-05 14:35:54.275156"
    54 = "2021-10-05 14:35:54.275159"
    55 = "2021-10-05 14:35:54.275161"
    56 = "2021-10-05 14:35:54.275163"
    57 = "2021-10-05 14:35:54.275166"
    58 = "2021-10-05 14:35:54.275168"
    59 = "2021-10-05 14:35:54.275170"
    60 = "2021-10-05 14:35:54.275173"
    61 = "2021-10-05 14:35:54.275175"
    62 = "2021-10-05 14:35:54.275177"
    63 = "2021-10-05 14:35:54.275179"
    64 = "2021-10-05 14:35:54.275182"
    65 = "2021-10-05 14:35:54.275184"
    66 = "2021-10-05 14:35:54.275187"
    67 = "2021-10-05 14:35:54.275190"
    68 = "2021-10-05 14:35:54.275192"
    69 = "2021-10-05 14:35:54.275195"
    70 = "2021-10-05 14:35:54.275197"
    71 = "2021-10-05 14:35:54.275199"
    72 = "2021-10-05 14:35:54.275202"
    73 = "2021-10-05 14:35:54.275204"
    74 = "2021-10-05 14:35:54.275207"
    75 = "2021-10-05 14:35:54.275209"
    76 = "2021-10-05 14:35:54.275212"
    77 = "2021-10-05 14:35:54.275214"
    78 = "2021-10-05 14:35:54.275218"
    79 = "2021-10-05 14:35:54.275221"
    80 = "2021-10-05 14:35:54.275224"
    81 = "2021-10-05 14:35:54.275226"
    82 = "2021-10-05 14:35:54.275228"
    83 = "2021-10-05 14:35:54.275231"
    84 = "2021-10-05 14:35:54.275233"
    85 = "2021-10-05 14:35:54.275235"
    86 = "2021-10-05 14:35:54.275238"
    87 = "2021-10-05 14:35:54.275240"
    88 = "2021-10-05 14:35:54.275243"
    89 = "2021-10-05 14:35:54.275245"
    90 = "2021-10-05 14:35:54.275247"
    91 = "2021-10-05 14:35:54.275250"
    92 = "2021-10-05 14:35:54.275252"
    93 = "2021-10-05 14:35:54.275254"
    94 = "2021-10-05 14:35:54.275258"
    95 = "2021-10-05 14:35:54.275260"
    96 = "2021-10-05 14:35:54.275265"
    97 = "2021-10-05 14:35:54.275267"
    98 = "2021-10-05 14:35:54.275270"
    99 = "2021-10-05 14:35:54.275272"
    100 = "2021-10-05 14:35:54.275275"
    101 = "2021-10-05 14:35:54.275277"
    102 = "2021-10-05 14:35:54.275279"
    103 = "2021-10-05 14:35:54.275282"
    104 = "2021-10-05 14:35:54.275284"
    105 = "2021-10-05 14:35:54.275286"
    106 = "2021-10-05 14:35:54.275289"
    107 = "2021-10-05 14:35:54.275291"
    108 = "2021-10-05 14:35:54.275293"
    109 = "2021-10-05 14:35:54.275296"
    110 = "2021-10-05 14:35:54.275300"
    111 = "2021-10-05 14:35:54.275303"
    112 = "2021-10-05 14:35:54.275306"
    113 = "2021-10-05 14:35:54.275308"
    114 = "2021-10-05 14:35:54.275310"
    115 = "2021-10-05 14:35:54.275313"
    116 = "2021-10-05 14:35:54.275315"
    117 = "2021-10-05 14:35:54.275318"
}

foreach ($row in $timeTaken.Keys) {
    $dataSheet.Range("F$row").Value = $timeTaken[$row]
}

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" worksheet right after "data" with panel-level info.
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Vasculitis"
$metaSheet.Range("C2").Value = 32
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.35"
$metaSheet.Range("D2").ClearFormats()
$metaSheet.Range("E2").Value = "2021-08-17T08:08:21.901398Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:54.271927"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/32/?format=json"

# ---------------------------------------------------------------------------
# 3) Match the header styling used on the "data" sheet (bold, bordered,
#    centered) for the new sheet's header row and its A2 index cell.
# ---------------------------------------------------------------------------
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$metaSheet.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Keep "data" as the active sheet/tab, matching the source workbook's
# unchanged <bookViews> (only the <sheets> list gained the new tab).
$dataSheet.Select()
$dataSheet.Range("A1").Select()

Write-Output "metadata sheet added; $($timeTaken.Count) timestamps refreshed"
